# C5-PowerPoint.pptx edit
#
# 1) Slide 6 ("SOURCES OF FINANCE") has a 4-column table whose table style
#    was changed to {4DE88A53-96A3-49D6-B7B6-3DE366B582A5}.
# 2) The deck's two theme parts had their colour palettes swapped: the
#    slide-master theme ("Integral") picks up the plain "Office Theme"
#    palette, while the notes-master theme becomes "Integral".
#    (The notes-master theme is not reachable through the Presentation
#    object model, so only the slide-master side of the swap -- the part
#    that actually affects what is shown on the slides -- can be applied
#    here; it is done through the per-slide ThemeColorScheme, which is the
#    supported way to repaint the 12 theme colour slots.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ---------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{4DE88A53-96A3-49D6-B7B6-3DE366B582A5}")

# --- 2. Swap the active theme's colour scheme over to "Office Theme" --
# Order matches the OOXML <a:clrScheme> child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB COM values are packed as 0xBBGGRR.
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Colors($i).RGB = $bgr
}
